# Auto-generated edit script applying the Marilith_Profits.xlsx diff
# to the corresponding sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# of the workbook. Each hunk in the source diff corresponds to one leve
# row whose price/profit columns (H..N) were refreshed by the scheduled
# market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 862.3
$ws.Range("I2").Value = 152.625
$ws.Range("J2").Value = 3701
$ws.Range("K2").Value = 152.625
$ws.Range("L2").Value = 3701
$ws.Range("M2").Value = -39.625
$ws.Range("N2").Value = -3927
$ws.Range("H28").Value = 212
$ws.Range("I28").Value = 203.33333
$ws.Range("J28").Value = 225
$ws.Range("K28").Value = 203.33333
$ws.Range("L28").Value = 225
$ws.Range("M28").Value = 281.66667
$ws.Range("N28").Value = -1195
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42496
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -132480
$ws.Range("H111").Value = 229
$ws.Range("I111").Value = 229
$ws.Range("K111").Value = 687
$ws.Range("M111").Value = 2380
$ws.Range("H115").Value = 246
$ws.Range("I115").Value = 246
$ws.Range("K115").Value = 738
$ws.Range("M115").Value = 829
$ws.Range("H125").Value = 254661.5
$ws.Range("I125").Value = 3374
$ws.Range("K125").Value = 30366
$ws.Range("M125").Value = -27906
$ws.Range("H135").Value = 975
$ws.Range("I135").Value = 975
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8775
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6240
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 831.125
$ws.Range("I74").Value = 879.63635
$ws.Range("K74").Value = 879.63635
$ws.Range("M74").Value = -5.636349999999993
$ws.Range("H77").Value = 831.125
$ws.Range("I77").Value = 879.63635
$ws.Range("K77").Value = 4398.18175
$ws.Range("M77").Value = -30.18174999999974
$ws.Range("H122").Value = 2077
$ws.Range("I122").Value = 2232.1667
$ws.Range("J122").Value = 1766.6666
$ws.Range("K122").Value = 6696.500100000001
$ws.Range("L122").Value = 5299.9998
$ws.Range("M122").Value = -4246.500100000001
$ws.Range("N122").Value = -10199.9998
$ws.Range("H123").Value = 57249.75
$ws.Range("J123").Value = 57249.75
$ws.Range("L123").Value = 57249.75
$ws.Range("N123").Value = -67049.75
$ws.Range("H132").Value = 3046.7827
$ws.Range("I132").Value = 2714.4736
$ws.Range("K132").Value = 8143.4208
$ws.Range("M132").Value = -5613.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 182
$ws.Range("I22").Value = 182
$ws.Range("K22").Value = 182
$ws.Range("M22").Value = -9
$ws.Range("H94").Value = 2240.7334
$ws.Range("I94").Value = 1622.9
$ws.Range("K94").Value = 1622.9
$ws.Range("M94").Value = -1171.9
$ws.Range("H105").Value = 4937.8
$ws.Range("I105").Value = 4937.8
$ws.Range("K105").Value = 4937.8
$ws.Range("M105").Value = -3190.8
$ws.Range("H107").Value = 1165
$ws.Range("I107").Value = 998
$ws.Range("K107").Value = 998
$ws.Range("M107").Value = 922
$ws.Range("H134").Value = 8676.637000000001
$ws.Range("I134").Value = 8194.333000000001
$ws.Range("K134").Value = 24582.999
$ws.Range("M134").Value = -22047.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 2000
$ws.Range("K11").Value = 2000
$ws.Range("M11").Value = -1860
$ws.Range("H31").Value = 1190.6875
$ws.Range("I31").Value = 1191.6666
$ws.Range("K31").Value = 1191.6666
$ws.Range("M31").Value = -896.6666
$ws.Range("H34").Value = 1190.6875
$ws.Range("I34").Value = 1191.6666
$ws.Range("K34").Value = 1191.6666
$ws.Range("M34").Value = -989.6666
$ws.Range("H58").Value = 1755.1111
$ws.Range("I58").Value = 1529.3572
$ws.Range("J58").Value = 1998.2307
$ws.Range("K58").Value = 1529.3572
$ws.Range("L58").Value = 1998.2307
$ws.Range("M58").Value = -1326.3572
$ws.Range("N58").Value = -2404.2307
$ws.Range("H107").Value = 580.7368
$ws.Range("I107").Value = 522.8
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 522.8
$ws.Range("L107").Value = 798
$ws.Range("M107").Value = 1397.2
$ws.Range("N107").Value = -4638
$ws.Range("H122").Value = 1350
$ws.Range("I122").Value = 1350
$ws.Range("K122").Value = 4050
$ws.Range("M122").Value = -1600
$ws.Range("H134").Value = 3199.6287
$ws.Range("I134").Value = 1250.4286
$ws.Range("J134").Value = 4499.095
$ws.Range("K134").Value = 3751.2858
$ws.Range("L134").Value = 13497.285
$ws.Range("M134").Value = -1216.2858
$ws.Range("N134").Value = -18567.285
$ws.Range("H136").Value = 1755.1111
$ws.Range("I136").Value = 1529.3572
$ws.Range("J136").Value = 1998.2307
$ws.Range("K136").Value = 4588.071599999999
$ws.Range("L136").Value = 5994.6921
$ws.Range("M136").Value = -2038.071599999999
$ws.Range("N136").Value = -11094.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 695.1818
$ws.Range("I23").Value = 489.8
$ws.Range("J23").Value = 866.3333
$ws.Range("K23").Value = 1469.4
$ws.Range("L23").Value = 2598.9999
$ws.Range("M23").Value = -1234.4
$ws.Range("N23").Value = -3068.9999
$ws.Range("H68").Value = 1313.5
$ws.Range("I68").Value = 1001.6
$ws.Range("K68").Value = 3004.8
$ws.Range("M68").Value = -2193.8
$ws.Range("H71").Value = 1313.5
$ws.Range("I71").Value = 1001.6
$ws.Range("K71").Value = 9014.4
$ws.Range("M71").Value = -4958.4
$ws.Range("H80").Value = 7833.8335
$ws.Range("J80").Value = 8834.333000000001
$ws.Range("L80").Value = 26502.999
$ws.Range("N80").Value = -28374.999
$ws.Range("H83").Value = 7833.8335
$ws.Range("J83").Value = 8834.333000000001
$ws.Range("L83").Value = 79508.997
$ws.Range("N83").Value = -88868.997
$ws.Range("H113").Value = 732.9231
$ws.Range("J113").Value = 1358
$ws.Range("L113").Value = 4074
$ws.Range("N113").Value = -8414
$ws.Range("H134").Value = 251024.5
$ws.Range("I134").Value = 251024.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 753073.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -748003.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 372642.84
$ws.Range("I3").Value = 400000
$ws.Range("J3").Value = 336166.66
$ws.Range("K3").Value = 400000
$ws.Range("L3").Value = 336166.66
$ws.Range("M3").Value = -399884
$ws.Range("N3").Value = -336398.66
$ws.Range("H10").Value = 6750
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 6750
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 6750
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -7088
$ws.Range("H11").Value = 746866.8
$ws.Range("I11").Value = 375012.88
$ws.Range("K11").Value = 375012.88
$ws.Range("M11").Value = -374873.88
$ws.Range("H70").Value = 2833333
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 2833333
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H99").Value = 15162
$ws.Range("I99").Value = 13194.4
$ws.Range("K99").Value = 13194.4
$ws.Range("M99").Value = -10948.4
$ws.Range("H102").Value = 12271.909
$ws.Range("I102").Value = 2487.6
$ws.Range("K102").Value = 2487.6
$ws.Range("M102").Value = -865.5999999999999
$ws.Range("H107").Value = 467.75
$ws.Range("I107").Value = 449
$ws.Range("K107").Value = 449
$ws.Range("M107").Value = 1471
$ws.Range("H126").Value = 9525.4
$ws.Range("I126").Value = 10666.333
$ws.Range("J126").Value = 7814
$ws.Range("K126").Value = 31998.999
$ws.Range("L126").Value = 23442
$ws.Range("M126").Value = -29528.999
$ws.Range("N126").Value = -28382
$ws.Range("H134").Value = 97059.28999999999
$ws.Range("J134").Value = 97059.28999999999
$ws.Range("L134").Value = 291177.87
$ws.Range("N134").Value = -296247.87

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6612.25
$ws.Range("I40").Value = 5350
$ws.Range("K40").Value = 5350
$ws.Range("M40").Value = -5214
$ws.Range("H122").Value = 3126.75
$ws.Range("I122").Value = 3072.7144
$ws.Range("K122").Value = 9218.143199999999
$ws.Range("M122").Value = -6768.143199999999
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 3875.8333
$ws.Range("I136").Value = 3501.8572
$ws.Range("J136").Value = 4399.4
$ws.Range("K136").Value = 10505.5716
$ws.Range("L136").Value = 13198.2
$ws.Range("M136").Value = -7955.571599999999
$ws.Range("N136").Value = -18298.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 34999.5
$ws.Range("J92").Value = 34999.5
$ws.Range("L92").Value = 34999.5
$ws.Range("N92").Value = -39991.5
$ws.Range("H122").Value = 3061.5833
$ws.Range("I122").Value = 3269.6
$ws.Range("J122").Value = 2913
$ws.Range("K122").Value = 9808.799999999999
$ws.Range("L122").Value = 8739
$ws.Range("M122").Value = -7358.799999999999
$ws.Range("N122").Value = -13639
$ws.Range("H126").Value = 31216.158
$ws.Range("I126").Value = 27394.889
$ws.Range("K126").Value = 82184.667
$ws.Range("M126").Value = -79714.667
$ws.Range("H132").Value = 2752.25
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 3001.6667
$ws.Range("K132").Value = 6012
$ws.Range("L132").Value = 9005.000100000001
$ws.Range("M132").Value = -3482
$ws.Range("N132").Value = -14065.0001

Write-Output "Applied 246 value updates and 6 clears across 8 sheets."
